$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the single placeholder email with the new list of three emails.
$ws.Range("A1").Value = "abibangbrandon855@gmail.com"
$ws.Range("A2").Value = "abibangbrandon86655@gmail.com"
$ws.Range("A3").Value = "abibangbrandon87755@gmail.com"

# The refreshed UI view is left-to-right (rightToLeft="0" on the sheet view).
try {
    $excel.ActiveWindow.DisplayRightToLeft = $false
} catch {
    # Not fatal if the host doesn't expose window-level view state.
}
